$d = $word.ActiveDocument

# Locate the abstract paragraph that needs its text rewritten and split into
# several runs (mirrors the original authoring history, including a
# proofErr-wrapped "analog" token) per the target diff. Find it by its
# distinctive original content rather than a hard-coded index, so the
# script is resilient to any unrelated structural differences.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*It was recently shown that reanalyses*") {
        $target = $p
        break
    }
}
if ($target -eq $null) {
    # Fallback: the abstract is the 3rd paragraph in the original layout.
    $target = $d.Paragraphs(3)
}
$r = $target.Range

$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="4EA3B5AD" w14:textId="7FB80D0E" w:rsidR="00793887" w:rsidRDefault="008F2EEF" w:rsidP="007B7786">' + `
    '<w:r w:rsidRPr="008F2EEF"><w:t>R</w:t></w:r>' + `
    '<w:r><w:t>eanalys</w:t></w:r>' + `
    '<w:r><w:t>i</w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve">s </w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve">datasets </w:t></w:r>' + `
    '<w:r><w:t>have an impact on statistical downscaling methods that may be even more important than the choice of</w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve"> the</w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve"> predictor variables. This work demonstrated the suitability of ERA5 over other global reanalyses for predicting daily precipitation at 301 stations in Switzerland, using six variants of </w:t></w:r>' + `
    '<w:proofErr w:type="spellStart"/>' + `
    '<w:r><w:t>analog</w:t></w:r>' + `
    '<w:proofErr w:type="spellEnd"/>' + `
    '<w:r><w:t xml:space="preserve"> methods. However, its high spatial resolution did not contribute to a gain in skill and was even counterproductive for simple calibration techniques.</w:t></w:r>' + `
    '</w:p>'

[void]$r.InsertXML($xml)
